# Fix the broken link to the GitHub repo on the "License" slide.
# The hyperlink text "Github repository" pointed at a stale/incorrect
# subfolder (left over from a different course's deck); it must point at
# this course's own folder in the teaching-material repo instead.

$p = $ppt.ActivePresentation

$oldAddress = "https://github.com/ETCE-LAB/teaching-material/tree/master/Emerging-Technologies-for-the-Circular-Economy"
$newAddress = "https://github.com/ETCE-LAB/teaching-material/tree/master/The-Limits-to-Growth"
$needle = "Github repository"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)

    # Identify the "License" slide by its title placeholder text.
    $titleShape = $slide.Shapes.Item(1)
    if (-not $titleShape.HasTextFrame) { continue }
    if ($titleShape.TextFrame.TextRange.Text -ne "License") { continue }

    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if (-not $shape.HasTextFrame) { continue }

        $tr = $shape.TextFrame.TextRange
        $full = $tr.Text
        if (-not $full) { continue }

        $idx = $full.IndexOf($needle)
        if ($idx -lt 0) { continue }

        $run = $tr.Characters($idx + 1, $needle.Length)
        $hyperlink = $run.ActionSettings.Item(1).Hyperlink

        if ($hyperlink.Address -eq $oldAddress) {
            $hyperlink.Address = $newAddress
        }
    }
}
